$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Delete empty "Title 26" placeholder shape
for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
  $shp = $s.Shapes.Item($i)
  if ($shp.Name -eq "Title 26") {
    $shp.Delete()
  }
}

$dx = 294928 / 12700.0
$dy = -1567296 / 12700.0
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
  $shp = $s.Shapes.Item($i)
  $shp.Left = $shp.Left + $dx
  $shp.Top = $shp.Top + $dy
}

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
  $shp = $s.Shapes.Item($i)
  Write-Output "$i : id=$($shp.Id) name=$($shp.Name) left=$($shp.Left) top=$($shp.Top)"
}
